# Fix property_category labels that were incorrectly left as "land" on the
# 建物 (Building) and 汽車 (Car) sheets.
#
# - 建物 (Building) sheet: column I is "property_category"; rows 2-5 were
#   all mistakenly tagged "land" and should read "building".
# - 汽車 (Car) sheet: column H is "property_category"; row 2 was mistakenly
#   tagged "land" and should read "car".

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"
$wsBuilding.Range("I5").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
